# Refresh the coinranking.com crypto snapshot: updates the Price (D) and 1h
# Volume (E) columns for every row, and swaps the Polkadot / WrappedEther
# rows (B, C, D, E) to reflect their new rank order, as produced by the
# scheduled "Updated cryptos list ... with GitHub Actions" job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($Addr, $Value)
    # The sheet stores Coin/Link/Price/Volume as literal text (inline strings),
    # even when a Price value such as "219.25" looks numeric. A leading quote
    # prefix keeps Excel from reinterpreting it as a number or date while
    # leaving the cell's number format untouched.
    $range = $ws.Range($Addr)
    if ($Value -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $range.Value = "'" + $Value
    } else {
        $range.Value = $Value
    }
}


Set-CellText "D2" "26.225.14"
Set-CellText "E2" "  -0.15%  "
Set-CellText "D3" "1.660.78"
Set-CellText "E3" "  -0.25%  "
Set-CellText "E4" "  -0.25%  "
Set-CellText "D5" "219.25"
Set-CellText "E5" "  +0.25%  "
Set-CellText "D6" "0.5259"
Set-CellText "E6" "  -0.82%  "
Set-CellText "E7" "  -0.28%  "
Set-CellText "D8" "0.2634"
Set-CellText "E8" "  -0.12%  "
Set-CellText "D9" "0.06319"
Set-CellText "E9" "  -0.72%  "
Set-CellText "D10" "20.63"
Set-CellText "E10" "  +0.30%  "
Set-CellText "D11" "0.07812"
Set-CellText "E11" "  -0.52%  "
Set-CellText "B12" "WrappedEther"
Set-CellText "C12" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-CellText "D12" "1.693.13"
Set-CellText "E12" "  +1.72%  "
Set-CellText "B13" "Polkadot"
Set-CellText "C13" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-CellText "D13" "4.491"
Set-CellText "E13" "  -1.64%  "
Set-CellText "D14" "1.889.72"
Set-CellText "E14" "  -0.15%  "
Set-CellText "D15" "0.5544"
Set-CellText "E15" "  +0.27%  "
Set-CellText "D16" "0.0₅8000"
Set-CellText "E16" "  -2.15%  "
Set-CellText "D17" "65.29"
Set-CellText "E17" "  -0.48%  "
Set-CellText "D18" "26.238.28"
Set-CellText "E18" "  -0.23%  "
Set-CellText "E19" "  -0.23%  "
Set-CellText "D20" "4.646"
Set-CellText "E20" "  -0.55%  "
Set-CellText "D21" "196.52"
Set-CellText "E21" "  +1.56%  "
Set-CellText "D22" "10.16"
Set-CellText "E22" "  -0.76%  "
Set-CellText "D23" "5.981"
Set-CellText "E24" "  -0.26%  "
Set-CellText "D25" "145.90"
Set-CellText "E25" "  +0.89%  "
Set-CellText "D26" "0.1204"
Set-CellText "E26" "  -1.82%  "
Set-CellText "D27" "7.165"
Set-CellText "E27" "  -0.71%  "
Set-CellText "D28" "16.06"
Set-CellText "E28" "  -0.20%  "
Set-CellText "D29" "1.518"
Set-CellText "E29" "  +2.41%  "
Set-CellText "D30" "0.05753"
Set-CellText "E30" "  -4.20%  "
Set-CellText "E31" "  -0.03%  "
Set-CellText "D32" "3.481"
Set-CellText "E32" "  -2.93%  "
Set-CellText "D33" "3.329"
Set-CellText "E33" "  +1.24%  "
Set-CellText "D34" "1.580"
Set-CellText "D35" "2.813"
Set-CellText "E35" "  -0.44%  "
Set-CellText "D36" "0.9523"
Set-CellText "E36" "  -0.93%  "
Set-CellText "D37" "2.427"
Set-CellText "E37" "  +0.05%  "
Set-CellText "D38" "0.5755"
Set-CellText "E38" "  -1.04%  "
Set-CellText "D39" "0.01592"
Set-CellText "E39" "  -0.86%  "
Set-CellText "D40" "5.948"
Set-CellText "E40" "  +1.24%  "
Set-CellText "D41" "1.058.37"
Set-CellText "E41" "  +0.82%  "
Set-CellText "D42" "0.8539"
Set-CellText "E42" "  -1.45%  "
Set-CellText "E43" "  -0.26%  "
Set-CellText "D44" "102.85"
Set-CellText "E44" "  -1.40%  "
Set-CellText "D45" "1.801.81"
Set-CellText "E45" "  -0.18%  "
Set-CellText "D46" "58.22"
Set-CellText "E46" "  +1.42%  "
Set-CellText "D47" "1.008"
Set-CellText "E47" "  -0.84%  "
Set-CellText "D48" "0.4415"
Set-CellText "E48" "  +0.77%  "
Set-CellText "D49" "7.995"
Set-CellText "E49" "  +0.02%  "
Set-CellText "D50" "0.05197"
Set-CellText "E50" "  +0.60%  "
Set-CellText "E51" "  -4.84%  "
